# Apply updated "dSF" (column F) values per the repull/recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    9  = 0
    11 = -2
    14 = -2
    17 = -13
    18 = -5
    22 = -3
    23 = 8
    25 = -3
    35 = -5
    52 = -2
    53 = -1
    59 = -2
    63 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
